{"js": "// Office.js (Word JavaScript API) script\n// Replicates the commit diff: updates the title date and every value in\n// the 20x5 arithmetic-expression table. Old cell texts are not all\n// globally unique after earlier replacements are applied (e.g. a later\n// cell's original text can equal an earlier cell's new text), so each\n// cell is addressed by its (row, col) position via Table.getCell and the\n// search/replace is scoped to that single cell's body. This keeps the\n// mapping unambiguous and preserves each run's formatting (font/size) and\n// paragraph alignment because insertText(..., Replace) only swaps the\n// matched range's text.\n\n// 1) Title paragraph date/weekday update.\nconst titleOld = \"2024-10-07 Monday\";\nconst titleNew = \"2024-10-08 Tuesday\";\nconst titleResults = context.document.body.search(titleOld, { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length === 0) {\n  throw new Error(\"Title text not found: \" + titleOld);\n}\ntitleResults.items[0].insertText(titleNew, Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Table cell values, addressed by (row, col).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellUpdates = [\n  [0, 0, \"81-24=\", \"7+46=\"],\n  [0, 1, \"60-35=\", \"19+4=\"],\n  [0, 2, \"54+18=\", \"74-19=\"],\n  [0, 3, \"40-3=\", \"73-66=\"],\n  [0, 4, \"91-42=\", \"96-68=\"],\n  [1, 0, \"16+27=\", \"58+27=\"],\n  [1, 1, \"93-86=\", \"97-18=\"],\n  [1, 2, \"35-28=\", \"92-35=\"],\n  [1, 3, \"36+19=\", \"41-12=\"],\n  [1, 4, \"53-4=\", \"41-39=\"],\n  [2, 0, \"70-27=\", \"34+39=\"],\n  [2, 1, \"8+26=\", \"48+46=\"],\n  [2, 2, \"74-69=\", \"6+26=\"],\n  [2, 3, \"88+6=\", \"90-31=\"],\n  [2, 4, \"72-34=\", \"18+63=\"],\n  [3, 0, \"6+35=\", \"88+3=\"],\n  [3, 1, \"72-7=\", \"94-29=\"],\n  [3, 2, \"56-49=\", \"95-56=\"],\n  [3, 3, \"25+57=\", \"55-46=\"],\n  [3, 4, \"91-45=\", \"33+38=\"],\n  [4, 0, \"55-36=\", \"58-9=\"],\n  [4, 1, \"96-39=\", \"24+7=\"],\n  [4, 2, \"84-77=\", \"65-57=\"],\n  [4, 3, \"56+35=\", \"37+7=\"],\n  [4, 4, \"76+19=\", \"74-58=\"],\n  [5, 0, \"71-45=\", \"55-16=\"],\n  [5, 1, \"47-8=\", \"44+49=\"],\n  [5, 2, \"90-61=\", \"92-25=\"],\n  [5, 3, \"18+27=\", \"92-39=\"],\n  [5, 4, \"46+49=\", \"7+44=\"],\n  [6, 0, \"32-4=\", \"59+12=\"],\n  [6, 1, \"84-67=\", \"55-26=\"],\n  [6, 2, \"64+28=\", \"49+2=\"],\n  [6, 3, \"18+4=\", \"37+58=\"],\n  [6, 4, \"85-58=\", \"34+39=\"],\n  [7, 0, \"6+78=\", \"11-3=\"],\n  [7, 1, \"78+4=\", \"24-6=\"],\n  [7, 2, \"18+47=\", \"40-21=\"],\n  [7, 3, \"84-28=\", \"59+25=\"],\n  [7, 4, \"80-41=\", \"59+6=\"],\n  [8, 0, \"33-9=\", \"25-19=\"],\n  [8, 1, \"44+29=\", \"24+39=\"],\n  [8, 2, \"52-34=\", \"47+29=\"],\n  [8, 3, \"90-83=\", \"9+14=\"],\n  [8, 4, \"26-18=\", \"53-16=\"],\n  [9, 0, \"94-45=\", \"70-44=\"],\n  [9, 1, \"14+49=\", \"38+26=\"],\n  [9, 2, \"48+35=\", \"49+4=\"],\n  [9, 3, \"29+69=\", \"33-4=\"],\n  [9, 4, \"64-9=\", \"6+26=\"],\n  [10, 0, \"95-29=\", \"86+5=\"],\n  [10, 1, \"60-33=\", \"93-18=\"],\n  [10, 2, \"81-27=\", \"4+39=\"],\n  [10, 3, \"37+16=\", \"87+8=\"],\n  [10, 4, \"36-29=\", \"68+4=\"],\n  [11, 0, \"93-66=\", \"54+27=\"],\n  [11, 1, \"25+58=\", \"27+65=\"],\n  [11, 2, \"94-47=\", \"19+45=\"],\n  [11, 3, \"3+88=\", \"63+8=\"],\n  [11, 4, \"66+16=\", \"14+79=\"],\n  [12, 0, \"39+59=\", \"90-43=\"],\n  [12, 1, \"39+3=\", \"26+36=\"],\n  [12, 2, \"4+59=\", \"8+59=\"],\n  [12, 3, \"61-7=\", \"7+45=\"],\n  [12, 4, \"62-59=\", \"30-2=\"],\n  [13, 0, \"26-9=\", \"34-26=\"],\n  [13, 1, \"91-27=\", \"56-47=\"],\n  [13, 2, \"56-37=\", \"66+8=\"],\n  [13, 3, \"38+29=\", \"14+48=\"],\n  [13, 4, \"41-34=\", \"86+7=\"],\n  [14, 0, \"62-19=\", \"26+8=\"],\n  [14, 1, \"29+63=\", \"23+38=\"],\n  [14, 2, \"85-49=\", \"7+87=\"],\n  [14, 3, \"70-21=\", \"8+83=\"],\n  [14, 4, \"38+59=\", \"62-8=\"],\n  [15, 0, \"67+29=\", \"28+8=\"],\n  [15, 1, \"8+54=\", \"85-77=\"],\n  [15, 2, \"97-58=\", \"7+66=\"],\n  [15, 3, \"37-8=\", \"92-37=\"],\n  [15, 4, \"9+76=\", \"56+18=\"],\n  [16, 0, \"41-33=\", \"91-45=\"],\n  [16, 1, \"19+35=\", \"9+25=\"],\n  [16, 2, \"94-88=\", \"65-26=\"],\n  [16, 3, \"71-23=\", \"59+13=\"],\n  [16, 4, \"71-65=\", \"27+27=\"],\n  [17, 0, \"55-29=\", \"64-46=\"],\n  [17, 1, \"81-49=\", \"46-7=\"],\n  [17, 2, \"80-9=\", \"76-28=\"],\n  [17, 3, \"95-76=\", \"84-45=\"],\n  [17, 4, \"34-26=\", \"10-5=\"],\n  [18, 0, \"29+9=\", \"61-28=\"],\n  [18, 1, \"25+59=\", \"19+4=\"],\n  [18, 2, \"54+29=\", \"29+2=\"],\n  [18, 3, \"34+49=\", \"80-1=\"],\n  [18, 4, \"69+3=\", \"36-27=\"],\n  [19, 0, \"45+38=\", \"79+7=\"],\n  [19, 1, \"56-17=\", \"90-55=\"],\n  [19, 2, \"45-29=\", \"15+19=\"],\n  [19, 3, \"17+14=\", \"50-43=\"],\n  [19, 4, \"65+16=\", \"41-35=\"],\n];\n\nfor (const [row, col, oldText, newText] of cellUpdates) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Cell (${row}, ${col}) text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Replicates the commit diff: updates the title date/weekday and every\n# value in the 20x5 arithmetic-expression table.\n#\n# Each cell is addressed by its (row, col) position via Table.Cell(row,\n# col) rather than by searching for its old text. This is deliberate:\n# some of the replacement values collide with other cells' original (or\n# freshly written) text \u2014 e.g. one cell's old text (\"26-9=\") becomes\n# another cell's new text (\"34-26=\"), and that exact string (\"34-26=\")\n# is also the *original* text of a later cell in the table. A text-based\n# Find/Replace (even one scoped to a single cell's Range) can therefore\n# land on the wrong cell once earlier edits have introduced a duplicate.\n# Writing directly to each cell's Range by position sidesteps that\n# ambiguity entirely and is safe regardless of processing order.\n#\n# The trailing character of a paragraph/cell Range.Text is the\n# paragraph mark (cells end in a cell mark); Range.End - 1 excludes it so\n# only the visible text is replaced, which keeps the existing run's\n# formatting (font/size) and paragraph alignment intact.\n\n$d = $word.ActiveDocument\n\nfunction Set-RangeText($range, $newText) {\n    $textRange = $d.Range($range.Start, $range.End - 1)\n    $textRange.Text = $newText\n}\n\n# 1) Title paragraph date/weekday update (first paragraph, above the table).\n$titlePara = $d.Paragraphs.Item(1)\nSet-RangeText $titlePara.Range \"2024-10-08 Tuesday\"\n\n# 2) Table cell values, addressed by (row, col) \u2014 1-based like Table.Cell().\n$table = $d.Tables.Item(1)\n\n$cellUpdates = @(\n    @{ Row = 1; Col = 1; OldText = \"81-24=\"; NewText = \"7+46=\" },\n    @{ Row = 1; Col = 2; OldText = \"60-35=\"; NewText = \"19+4=\" },\n    @{ Row = 1; Col = 3; OldText = \"54+18=\"; NewText = \"74-19=\" },\n    @{ Row = 1; Col = 4; OldText = \"40-3=\"; NewText = \"73-66=\" },\n    @{ Row = 1; Col = 5; OldText = \"91-42=\"; NewText = \"96-68=\" },\n    @{ Row = 2; Col = 1; OldText = \"16+27=\"; NewText = \"58+27=\" },\n    @{ Row = 2; Col = 2; OldText = \"93-86=\"; NewText = \"97-18=\" },\n    @{ Row = 2; Col = 3; OldText = \"35-28=\"; NewText = \"92-35=\" },\n    @{ Row = 2; Col = 4; OldText = \"36+19=\"; NewText = \"41-12=\" },\n    @{ Row = 2; Col = 5; OldText = \"53-4=\"; NewText = \"41-39=\" },\n    @{ Row = 3; Col = 1; OldText = \"70-27=\"; NewText = \"34+39=\" },\n    @{ Row = 3; Col = 2; OldText = \"8+26=\"; NewText = \"48+46=\" },\n    @{ Row = 3; Col = 3; OldText = \"74-69=\"; NewText = \"6+26=\" },\n    @{ Row = 3; Col = 4; OldText = \"88+6=\"; NewText = \"90-31=\" },\n    @{ Row = 3; Col = 5; OldText = \"72-34=\"; NewText = \"18+63=\" },\n    @{ Row = 4; Col = 1; OldText = \"6+35=\"; NewText = \"88+3=\" },\n    @{ Row = 4; Col = 2; OldText = \"72-7=\"; NewText = \"94-29=\" },\n    @{ Row = 4; Col = 3; OldText = \"56-49=\"; NewText = \"95-56=\" },\n    @{ Row = 4; Col = 4; OldText = \"25+57=\"; NewText = \"55-46=\" },\n    @{ Row = 4; Col = 5; OldText = \"91-45=\"; NewText = \"33+38=\" },\n    @{ Row = 5; Col = 1; OldText = \"55-36=\"; NewText = \"58-9=\" },\n    @{ Row = 5; Col = 2; OldText = \"96-39=\"; NewText = \"24+7=\" },\n    @{ Row = 5; Col = 3; OldText = \"84-77=\"; NewText = \"65-57=\" },\n    @{ Row = 5; Col = 4; OldText = \"56+35=\"; NewText = \"37+7=\" },\n    @{ Row = 5; Col = 5; OldText = \"76+19=\"; NewText = \"74-58=\" },\n    @{ Row = 6; Col = 1; OldText = \"71-45=\"; NewText = \"55-16=\" },\n    @{ Row = 6; Col = 2; OldText = \"47-8=\"; NewText = \"44+49=\" },\n    @{ Row = 6; Col = 3; OldText = \"90-61=\"; NewText = \"92-25=\" },\n    @{ Row = 6; Col = 4; OldText = \"18+27=\"; NewText = \"92-39=\" },\n    @{ Row = 6; Col = 5; OldText = \"46+49=\"; NewText = \"7+44=\" },\n    @{ Row = 7; Col = 1; OldText = \"32-4=\"; NewText = \"59+12=\" },\n    @{ Row = 7; Col = 2; OldText = \"84-67=\"; NewText = \"55-26=\" },\n    @{ Row = 7; Col = 3; OldText = \"64+28=\"; NewText = \"49+2=\" },\n    @{ Row = 7; Col = 4; OldText = \"18+4=\"; NewText = \"37+58=\" },\n    @{ Row = 7; Col = 5; OldText = \"85-58=\"; NewText = \"34+39=\" },\n    @{ Row = 8; Col = 1; OldText = \"6+78=\"; NewText = \"11-3=\" },\n    @{ Row = 8; Col = 2; OldText = \"78+4=\"; NewText = \"24-6=\" },\n    @{ Row = 8; Col = 3; OldText = \"18+47=\"; NewText = \"40-21=\" },\n    @{ Row = 8; Col = 4; OldText = \"84-28=\"; NewText = \"59+25=\" },\n    @{ Row = 8; Col = 5; OldText = \"80-41=\"; NewText = \"59+6=\" },\n    @{ Row = 9; Col = 1; OldText = \"33-9=\"; NewText = \"25-19=\" },\n    @{ Row = 9; Col = 2; OldText = \"44+29=\"; NewText = \"24+39=\" },\n    @{ Row = 9; Col = 3; OldText = \"52-34=\"; NewText = \"47+29=\" },\n    @{ Row = 9; Col = 4; OldText = \"90-83=\"; NewText = \"9+14=\" },\n    @{ Row = 9; Col = 5; OldText = \"26-18=\"; NewText = \"53-16=\" },\n    @{ Row = 10; Col = 1; OldText = \"94-45=\"; NewText = \"70-44=\" },\n    @{ Row = 10; Col = 2; OldText = \"14+49=\"; NewText = \"38+26=\" },\n    @{ Row = 10; Col = 3; OldText = \"48+35=\"; NewText = \"49+4=\" },\n    @{ Row = 10; Col = 4; OldText = \"29+69=\"; NewText = \"33-4=\" },\n    @{ Row = 10; Col = 5; OldText = \"64-9=\"; NewText = \"6+26=\" },\n    @{ Row = 11; Col = 1; OldText = \"95-29=\"; NewText = \"86+5=\" },\n    @{ Row = 11; Col = 2; OldText = \"60-33=\"; NewText = \"93-18=\" },\n    @{ Row = 11; Col = 3; OldText = \"81-27=\"; NewText = \"4+39=\" },\n    @{ Row = 11; Col = 4; OldText = \"37+16=\"; NewText = \"87+8=\" },\n    @{ Row = 11; Col = 5; OldText = \"36-29=\"; NewText = \"68+4=\" },\n    @{ Row = 12; Col = 1; OldText = \"93-66=\"; NewText = \"54+27=\" },\n    @{ Row = 12; Col = 2; OldText = \"25+58=\"; NewText = \"27+65=\" },\n    @{ Row = 12; Col = 3; OldText = \"94-47=\"; NewText = \"19+45=\" },\n    @{ Row = 12; Col = 4; OldText = \"3+88=\"; NewText = \"63+8=\" },\n    @{ Row = 12; Col = 5; OldText = \"66+16=\"; NewText = \"14+79=\" },\n    @{ Row = 13; Col = 1; OldText = \"39+59=\"; NewText = \"90-43=\" },\n    @{ Row = 13; Col = 2; OldText = \"39+3=\"; NewText = \"26+36=\" },\n    @{ Row = 13; Col = 3; OldText = \"4+59=\"; NewText = \"8+59=\" },\n    @{ Row = 13; Col = 4; OldText = \"61-7=\"; NewText = \"7+45=\" },\n    @{ Row = 13; Col = 5; OldText = \"62-59=\"; NewText = \"30-2=\" },\n    @{ Row = 14; Col = 1; OldText = \"26-9=\"; NewText = \"34-26=\" },\n    @{ Row = 14; Col = 2; OldText = \"91-27=\"; NewText = \"56-47=\" },\n    @{ Row = 14; Col = 3; OldText = \"56-37=\"; NewText = \"66+8=\" },\n    @{ Row = 14; Col = 4; OldText = \"38+29=\"; NewText = \"14+48=\" },\n    @{ Row = 14; Col = 5; OldText = \"41-34=\"; NewText = \"86+7=\" },\n    @{ Row = 15; Col = 1; OldText = \"62-19=\"; NewText = \"26+8=\" },\n    @{ Row = 15; Col = 2; OldText = \"29+63=\"; NewText = \"23+38=\" },\n    @{ Row = 15; Col = 3; OldText = \"85-49=\"; NewText = \"7+87=\" },\n    @{ Row = 15; Col = 4; OldText = \"70-21=\"; NewText = \"8+83=\" },\n    @{ Row = 15; Col = 5; OldText = \"38+59=\"; NewText = \"62-8=\" },\n    @{ Row = 16; Col = 1; OldText = \"67+29=\"; NewText = \"28+8=\" },\n    @{ Row = 16; Col = 2; OldText = \"8+54=\"; NewText = \"85-77=\" },\n    @{ Row = 16; Col = 3; OldText = \"97-58=\"; NewText = \"7+66=\" },\n    @{ Row = 16; Col = 4; OldText = \"37-8=\"; NewText = \"92-37=\" },\n    @{ Row = 16; Col = 5; OldText = \"9+76=\"; NewText = \"56+18=\" },\n    @{ Row = 17; Col = 1; OldText = \"41-33=\"; NewText = \"91-45=\" },\n    @{ Row = 17; Col = 2; OldText = \"19+35=\"; NewText = \"9+25=\" },\n    @{ Row = 17; Col = 3; OldText = \"94-88=\"; NewText = \"65-26=\" },\n    @{ Row = 17; Col = 4; OldText = \"71-23=\"; NewText = \"59+13=\" },\n    @{ Row = 17; Col = 5; OldText = \"71-65=\"; NewText = \"27+27=\" },\n    @{ Row = 18; Col = 1; OldText = \"55-29=\"; NewText = \"64-46=\" },\n    @{ Row = 18; Col = 2; OldText = \"81-49=\"; NewText = \"46-7=\" },\n    @{ Row = 18; Col = 3; OldText = \"80-9=\"; NewText = \"76-28=\" },\n    @{ Row = 18; Col = 4; OldText = \"95-76=\"; NewText = \"84-45=\" },\n    @{ Row = 18; Col = 5; OldText = \"34-26=\"; NewText = \"10-5=\" },\n    @{ Row = 19; Col = 1; OldText = \"29+9=\"; NewText = \"61-28=\" },\n    @{ Row = 19; Col = 2; OldText = \"25+59=\"; NewText = \"19+4=\" },\n    @{ Row = 19; Col = 3; OldText = \"54+29=\"; NewText = \"29+2=\" },\n    @{ Row = 19; Col = 4; OldText = \"34+49=\"; NewText = \"80-1=\" },\n    @{ Row = 19; Col = 5; OldText = \"69+3=\"; NewText = \"36-27=\" },\n    @{ Row = 20; Col = 1; OldText = \"45+38=\"; NewText = \"79+7=\" },\n    @{ Row = 20; Col = 2; OldText = \"56-17=\"; NewText = \"90-55=\" },\n    @{ Row = 20; Col = 3; OldText = \"45-29=\"; NewText = \"15+19=\" },\n    @{ Row = 20; Col = 4; OldText = \"17+14=\"; NewText = \"50-43=\" },\n    @{ Row = 20; Col = 5; OldText = \"65+16=\"; NewText = \"41-35=\" }\n)\n\nforeach ($update in $cellUpdates) {\n    $cell = $table.Cell($update.Row, $update.Col)\n    Set-RangeText $cell.Range $update.NewText\n}\n"}
